# Readme v2.0.0 Updated / 기능명세 Updated
# Updates the "클래스"(Class) / "속성"(Property) reference columns (E/F)
# on the "기능" sheet to reflect the renamed internal data-model symbols.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("기능")

# Row 2: 전체 CPU 사용률
$ws.Cells.Item(2, 5).Value = "PerfDataOSProcessor"
$ws.Cells.Item(2, 6).Value = "m_table.cpuTotal"

# Row 3: 프로세스 이름 (CPU 섹션)
$ws.Cells.Item(3, 5).Value = "CPerfDataPerProcess"
$ws.Cells.Item(3, 6).Value = "PerProcessDataObj.name"

# Row 4: PID (CPU 섹션)
$ws.Cells.Item(4, 5).Value = "CPerfDataPerProcess"
$ws.Cells.Item(4, 6).Value = "m_table.first"

# Row 5: 프로세스 CPU 평균(60초) 사용률
$ws.Cells.Item(5, 5).Value = "CPerfDataPerProcess"
$ws.Cells.Item(5, 6).Value = "PerProcessDataObj.meanUsageRate"

# Row 6: 프로세스 CPU 사용률
$ws.Cells.Item(6, 5).Value = "CPerfDataPerProcess"
$ws.Cells.Item(6, 6).Value = "PerProcessDataObj.usageRate"

# Row 7: 전체 메모리
$ws.Cells.Item(7, 5).Value = "CPerfDataOS"
$ws.Cells.Item(7, 6).Value = "OSDataObj.totalVisibleMemory"

# Row 8: 전체 사용중인 메모리
# (leading apostrophe preserves the original "quote prefix" text-literal cell
# style - s="7" / quotePrefix="1" - that this cell already carried)
$ws.Cells.Item(8, 5).Value = "CPerfDataOS"
$ws.Cells.Item(8, 6).Value = "'OSDataObj.freePhysicalMemory"

# Row 9: 프로세스 이름 (메모리 섹션)
$ws.Cells.Item(9, 5).Value = "CPerfDataPerProcess"
$ws.Cells.Item(9, 6).Value = "PerProcessDataObj.name"

# Row 10: PID (메모리 섹션)
$ws.Cells.Item(10, 5).Value = "CPerfDataPerProcess"
$ws.Cells.Item(10, 6).Value = "m_table.first"

# Row 11: 프로세스 커밋 메모리
$ws.Cells.Item(11, 5).Value = "CPerfDataPerProcess"
$ws.Cells.Item(11, 6).Value = "PerProcessDataObj.virtualBytes"

# Row 12: 프로세스 개인 메모리
$ws.Cells.Item(12, 5).Value = "CPerfDataPerProcess"
$ws.Cells.Item(12, 6).Value = "PerProcessDataObj.privateBytes"

# Row 13: 프로세스 작업집합
$ws.Cells.Item(13, 5).Value = "CPerfDataPerProcess"
$ws.Cells.Item(13, 6).Value = "PerProcessDataObj.workingSet"

# Row 14: 디스크(파티션)이름
$ws.Cells.Item(14, 6).Value = "LogicalDiskDataObj.deviceID"

# Row 15: 사용 중 용량
$ws.Cells.Item(15, 6).Value = "LogicalDiskDataObj.size"

# Row 16: 사용 가능 용량
$ws.Cells.Item(16, 6).Value = "LogicalDiskDataObj.freeSpace"

# Row 17: 프로세스 이름 (디스크 프로세스별 사용량 섹션)
$ws.Cells.Item(17, 5).Value = "CPerfDataPerProcess"
$ws.Cells.Item(17, 6).Value = "PerProcessDataObj.name"

# Row 18: PID (디스크 프로세스별 사용량 섹션)
$ws.Cells.Item(18, 5).Value = "Etw"
$ws.Cells.Item(18, 6).Value = "diskMap.first"

# Row 19: 프로세스 읽기 속도
$ws.Cells.Item(19, 5).Value = "Etw"
$ws.Cells.Item(19, 6).Value = "ProcessDiskData.readBytes"

# Row 20: 프로세스 쓰기 속도
$ws.Cells.Item(20, 5).Value = "Etw"
$ws.Cells.Item(20, 6).Value = "ProcessDiskData.writeBytes"

# Row 21: 전체 네트워크 I/O속도
$ws.Cells.Item(21, 5).Value = "Etw"
$ws.Cells.Item(21, 6).Value = "networkMap 추가 연산"

# Row 22: 프로세스 이름 (네트워크 프로세스별 사용률 섹션)
$ws.Cells.Item(22, 5).Value = "CPerfDataPerProcess"
$ws.Cells.Item(22, 6).Value = "Name"

# Row 23: PID (네트워크 프로세스별 사용률 섹션)
$ws.Cells.Item(23, 5).Value = "Etw"
$ws.Cells.Item(23, 6).Value = "networkMap.first"

# Row 24: 보내기 속도
$ws.Cells.Item(24, 5).Value = "Etw"
$ws.Cells.Item(24, 6).Value = "ProcessNetworkData.sendBytes"

# Row 25: 받기 속도
$ws.Cells.Item(25, 5).Value = "Etw"
$ws.Cells.Item(25, 6).Value = "ProcessNetworkData.receiveBytes"
